$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target final values for the rows whose player records were reshuffled.
$ws.Range("A4").Value = "Stephen Curry"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Golden State Warriors"

$ws.Range("A5").Value = "Tyrese Haliburton"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Indiana Pacers"

$ws.Range("A6").Value = "OG Anunoby"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "New York Knicks"

$ws.Range("A10").Value = "Daniel Gafford"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Dallas Mavericks"

$ws.Range("A14").Value = "Keegan Murray"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Sacramento Kings"

$ws.Range("A16").Value = "Jarrett Allen"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Cleveland Cavaliers"

$ws.Range("A19").Value = "Darius Garland"
$ws.Range("B19").Value = "PG"
$ws.Range("C19").Value = "Cleveland Cavaliers"

$wb.Save()
